# Add audit files to data library (#3365)
#
# The JASP example/data-library description sheet gains two new rows
# describing the "BuildIt" audit example files. These are inserted above
# the existing "Miscellaneous" block (which shifts down by two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the old row 50 ("Anscombe's Quartet"
# / Miscellaneous block), pushing everything from old row 50 onward down
# by two rows.
$ws.Rows("50:51").Insert()

# New row 50: BuildIt Binary (non-monetary audit population)
# New row 51: BuildIt Monetary (monetary audit population)
# Column order: D=Raw Data, E=Chapter, F=Analysis, G=Description, H=JASP
# file commented, I=Variable types correct (A/B/C/J stay blank for these
# two rows, matching the source data).
$ws.Range("F50").Value = "Audit"
$ws.Range("E50").Value = "BuildIt Binary"
$ws.Range("E51").Value = "BuildIt Monetary"
$ws.Range("G51").Value = "BuildIt's monetary audit population consisting of 3500 transactions."
$ws.Range("G50").Value = "BuildIt's non-monetary audit population consisting of 3500 records. "

$ws.Range("D50").Value = "Yes"
$ws.Range("H50").Value = "No"
$ws.Range("I50").Value = "Yes"

$ws.Range("D51").Value = "Yes"
$ws.Range("F51").Value = "Audit"
$ws.Range("H51").Value = "No"
$ws.Range("I51").Value = "Yes"

# New column of widths for the now-used column H, matching the rest of
# the description columns.
$ws.Columns.Item(8).ColumnWidth = 19

# Page setup was (re)established for this sheet when it was resaved.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Reflect the final cursor/selection position recorded in the workbook.
$ws.Range("I52").Select()

Write-Output "Audit rows inserted."
